# ---------------------------------------------------------------------------
# Sensitivity analysis for C and k: refreshed the "Calculated in situ DIC
# (umol/kg)" values in column T (older titrant-molinity script results).
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tValues = @{
    29 = 1728.4111633839241
    30 = 1939.0750753074369
    31 = 2126.0408110840708
    32 = 1852.235109647253
    33 = 1764.065435550898
    34 = 1640.6311129712799
    36 = 1702.772989296077
    37 = 1972.4841360953819
    38 = 2118.847393016958
    39 = 2078.7387755524251
    40 = 2116.5155141367659
    41 = 1992.1652563154039
    42 = 1982.8895791860141
    43 = 1977.83536937958
    44 = 1845.2477594969309
    45 = 1860.4433565011759
    46 = 1734.1488305150981
    47 = 1747.5027864090371
    56 = 1727.0774261756301
    57 = 1762.4220047633039
    58 = 1799.9086993448341
    59 = 1784.869886623304
    60 = 1773.5304510336
    61 = 1754.085190201361
    62 = 1770.142780273223
    63 = 1747.5409012868411
    64 = 1712.9509906108369
    65 = 1698.1483636858829
    66 = 1705.8300826620609
    67 = 1706.171345967027
    68 = 2235.7270752876379
    69 = 1708.0644987931701
    70 = 1689.3489301596831
    71 = 1720.1827637339111
    72 = 1714.6523915016569
    74 = 2076.5172146505452
    75 = 2080.8704689344481
    76 = 2075.658816176599
    77 = 2098.2084397638018
    78 = 2186.9378101158668
    79 = 2075.9298813237019
    80 = 2090.0033292014641
    81 = 2120.1700026324311
    82 = 2210.1353634709649
    83 = 1888.7220330848729
    84 = 2073.4999273761809
    85 = 2066.1122185424229
    86 = 1743.437571131338
    87 = 2139.3074571035509
    88 = 2133.8451916908102
    89 = 2109.1470511792209
    90 = 2134.1578721894898
    91 = 2138.3861010535688
    92 = 2170.180756087776
    93 = 2167.2432661556109
    102 = 1763.0486226437249
    103 = 2110.2873260929532
    104 = 2121.0763656656131
    105 = 2140.9961495842408
    106 = 2140.738026049884
    107 = 2262.8189410786222
    108 = 2157.4550962802318
    109 = 2227.9333775045179
    110 = 2207.050096666389
    111 = 2146.8455950201478
    114 = 2182.0897859422521
    115 = 2163.989167442759
    116 = 2188.0502296822938
    117 = 2187.496010443951
    118 = 2349.568037690859
    119 = 2208.469223095869
    120 = 2227.1294799252669
    121 = 2175.1137213648808
    122 = 2175.6419855208792
    123 = 2174.9428955531098
    124 = 2176.8750660763098
    125 = 2286.8264757218631
    126 = 2168.0859037459222
    127 = 1789.009232064225
    128 = 1784.258772400284
    130 = 1802.764062733839
    147 = 1735.7372861829399
    148 = 1772.8659065072859
    149 = 1880.822453594491
    150 = 1934.9097480683149
    151 = 2019.3270391437679
    152 = 2073.623426228211
    153 = 2118.2984501519868
    154 = 2140.623517699742
    155 = 2170.2135607426608
    156 = 2157.7062308030299
    157 = 2205.598343709356
    158 = 2242.81790761191
    159 = 2224.527965557711
    160 = 2045.002887320293
    161 = 2238.5442355774439
    174 = 1764.4177852175601
    176 = 1726.3147021265549
    177 = 1716.4041325756159
    178 = 1777.569657580184
    179 = 1839.455731279259
    180 = 1920.452988319531
    181 = 1994.5692524580711
    182 = 2038.6970139876371
    183 = 2080.8092191471778
    184 = 2113.5223767798489
    185 = 2157.3775729523691
    189 = 1704.091139768507
    190 = 1720.1443165635101
    191 = 1634.5514785489329
    192 = 1684.627716583296
    193 = 1734.4458451134369
    194 = 2171.5826220980239
    195 = 2208.905630646349
    196 = 1603.7542701999059
    198 = 1622.688928798196
    201 = 1713.418915056744
    202 = 1719.791343803839
    203 = 1989.821281995532
    204 = 1808.5989002578369
    205 = 1894.3752754033601
    206 = 1981.9899177264961
    207 = 2045.288324774782
    208 = 2208.1489202264911
    209 = 2077.9516864681141
    210 = 2130.482053894757
    211 = 2129.9832255722008
    212 = 1738.3624558021199
    213 = 2089.0630975833119
    214 = 2152.4406720651
    215 = 2161.2012511214061
    216 = 2170.022107750518
    217 = 2181.0795895811029
    218 = 2188.671506118847
    219 = 2195.6759000073789
    220 = 2201.4270390400789
    221 = 2202.3406802730292
    222 = 1764.44986759627
    225 = 1726.1277958392741
    226 = 2381.660400391173
    227 = 2092.304944573435
    228 = 2124.4802095075429
    229 = 2114.5287239757372
    230 = 2141.1632143872721
    231 = 2152.5605294461998
    232 = 2261.3745396517688
    233 = 2108.817360987257
    234 = 2114.832119163902
    235 = 2113.44539644919
    236 = 2104.9853041344941
    238 = 2193.692950369169
    239 = 2181.4821736323311
    240 = 2527.1411218953458
    241 = 2089.580488322877
    242 = 2127.4254589599468
    243 = 2119.7112863565881
    244 = 2117.3546768193469
    245 = 2135.4012085443419
    248 = 1711.2370372453661
    249 = 1724.174326100534
    250 = 1861.0396868001801
    251 = 1736.182162228429
    252 = 1745.5302759119991
    253 = 1770.5230988410381
    254 = 1743.3480425825551
    255 = 1729.5152185879119
    256 = 1708.716264853753
    257 = 1739.7260846603519
    258 = 1783.936474410365
    259 = 1947.304629626396
    260 = 1932.9910369035381
    261 = 1831.465713187971
    262 = 1933.758216992019
    263 = 1708.876557123416
    277 = 1701.196814155669
    278 = 1721.136943399674
    279 = 1737.670889194168
    280 = 1767.645544349394
    281 = 1786.226144057488
    282 = 1825.3668164913749
    283 = 1874.7850331129559
    284 = 1858.049517276547
    286 = 1874.8899007476559
    287 = 1934.021762050286
    288 = 1976.4005203837701
    289 = 2016.7719927577259
    290 = 1772.48999089173
    291 = 2045.0161418985069
    292 = 2085.8794748368819
    293 = 2115.8846816910718
    294 = 2137.6524488976111
    295 = 2163.788188520824
    296 = 2178.97144309021
    297 = 2188.6941409145052
    298 = 1718.50610711182
    299 = 2138.2909416625962
    301 = 1790.1023155204459
    302 = 1796.472659850243
    303 = 1797.823500374401
    304 = 1632.3432274745201
    305 = 1676.5490090633909
    306 = 1665.5669700790311
    307 = 1641.1717287585259
    308 = 1649.280500000239
    309 = 1717.1725237035801
    310 = 1706.427261984607
    311 = 1683.803017547692
    313 = 1619.7571697811161
    314 = 1851.4645413296489
    315 = 1785.2081992666881
    316 = 1774.6091137607641
    317 = 1806.4825229547951
    318 = 2200.438006317926
    319 = 1831.0295929996639
    320 = 1829.534427821425
    321 = 1858.000205621011
    322 = 1853.5901961346649
    324 = 1874.6132516016669
    325 = 1894.9444969068411
    326 = 1912.368193279769
    327 = 1936.4419918417859
    328 = 1966.1376140100101
    329 = 1983.761065272475
    331 = 2230.041489930738
    332 = 2021.2653315421021
    333 = 2056.0703879683729
    334 = 2101.0443275422522
    335 = 2128.68083839886
    336 = 2168.3609256734362
    337 = 2185.3190047344619
    338 = 2241.4753220648381
    339 = 1807.124441551008
}

foreach ($row in $tValues.Keys) {
    $ws.Cells.Item([int]$row, 20).Value = $tValues[$row]
}

# Restore the view roughly where the author left it: scrolled so column D is
# leftmost visible, with the whole of column T selected (active cell near row 6).
$ws.Range("D6").Select()
$ws.Columns("T:T").Select()

Write-Output "Updated $($tValues.Count) cells in column T"
